$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted as the new top row of the
# "Perejil" data block (row 661), pushing all subsequent rows down by one.
$ws.Rows("661:661").Insert()

$ws.Range("A661").Value = 6
$ws.Range("B661").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C661").Value = "Metropolitana"
$ws.Range("D661").Value = 45077
$ws.Range("D661").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E661").Value = 13
$ws.Range("F661").Value = 100112044
$ws.Range("G661").Value = "Perejil"
$ws.Range("H661").Value = "Sin especificar"
$ws.Range("I661").Value = "Primera"
$ws.Range("J661").Value = 250
$ws.Range("K661").Value = 11000
$ws.Range("L661").Value = 12000
$ws.Range("M661").Value = 11440
$ws.Range("N661").Value = "`$/docena de atados"
$ws.Range("O661").Value = "Región Metropolitana"
$ws.Range("P661").Value = 3813
$ws.Range("Q661").Value = 3
$ws.Range("R661").Value = "Hortaliza"
